# Summarize findings for google
# Updates the RQ1 (Cause of Flakiness), RQ2 (Fix for Flakiness), and
# Programming Language count tables on Sheet1 with refreshed totals and
# a handful of newly observed categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- RQ1: Cause of Flakiness? (column B/C) ----
$ws.Range("C4").Value = 24
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 20
$ws.Range("C8").Value = 6
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 3

# New category row: Bit Manipulation / Arithmetic
$ws.Range("B12").Value = "Bit Manipulation / Arithmetic"
$ws.Range("C12").Value = 1

# ---- RQ2: Fix for Flakiness? (column E/F) ----
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 7
$ws.Range("F9").Value = 4
$ws.Range("F13").Value = 7
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 2

# New fix rows
$ws.Range("E18").Value = "(Bit Manipulation / Arithmetic) bit clear"
$ws.Range("F18").Value = 1

$ws.Range("E19").Value = "(Concurrency) update global state"
$ws.Range("F19").Value = 1

# ---- Programming Language (column B/C) ----
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 8
$ws.Range("C23").Value = 5
$ws.Range("C25").Value = 17
$ws.Range("C27").Value = 7
$ws.Range("C28").Value = 7

# New language row: Rust
$ws.Range("B30").Value = "Rust"
$ws.Range("C30").Value = 1

$ws.Range("E20").Value = "(Async Wait) Promise statement"
$ws.Range("F20").Value = 1

# ---- Resize tables to include newly added rows ----
$wb.Worksheets.Item(1).ListObjects.Item("Table13").Resize($ws.Range("E3:F20"))
$wb.Worksheets.Item(1).ListObjects.Item("Table3").Resize($ws.Range("B18:C30"))

# ---- Update selection state to match the saved file ----
$ws.Range("C7").Select()
